$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BFDR")

# The sheet contains two fully-blank rows (132 and 136) that separate
# groups of related entries. Remove the line breaks by deleting these
# blank rows outright, which shifts the subsequent rows up.
# Delete the lower one first so the earlier row index is unaffected.
$ws.Rows.Item(136).Delete()
$ws.Rows.Item(132).Delete()

# Keep the hidden filter-database defined name in sync with the new
# (shorter) data range now that the two blank rows are gone.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "BFDR!_FilterDatabase") {
        $n.RefersTo = "=BFDR!`$A`$1:`$H`$139"
    }
}

# Refresh the sheet's remembered sort range so it reflects the two
# fewer rows as well.
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A139"))
$sort.SetRange($ws.Range("A2:H140"))
$sort.Header = 2
$sort.Apply()

# Leave the view focused on the BFDR sheet, near the rows that were
# just collapsed together.
$ws.Activate()
$ws.Range("A132:XFD132").Select()
